$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 46064

# Row 3
$ws.Range("C3").Value = 46064

# Row 4
$ws.Range("A4").Value = "A 53519-2023"
$ws.Range("B4").Value = 45230
$ws.Range("C4").Value = 46064
$ws.Range("G4").Value = 1
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1261/artfynd/A 53519-2023 artfynd.xlsx", "A 53519-2023")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1261/kartor/A 53519-2023 karta.png", "A 53519-2023")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1261/klagomål/A 53519-2023 FSC-klagomål.docx", "A 53519-2023")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1261/klagomålsmail/A 53519-2023 FSC-klagomål mail.docx", "A 53519-2023")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1261/tillsyn/A 53519-2023 tillsynsbegäran.docx", "A 53519-2023")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1261/tillsynsmail/A 53519-2023 tillsynsbegäran mail.docx", "A 53519-2023")'

# Row 5
$ws.Range("A5").Value = "A 50825-2025"
$ws.Range("B5").Value = 45946.54048611111
$ws.Range("C5").Value = 46064
$ws.Range("G5").Value = 2.6
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1261/artfynd/A 50825-2025 artfynd.xlsx", "A 50825-2025")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1261/kartor/A 50825-2025 karta.png", "A 50825-2025")'
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1261/klagomål/A 50825-2025 FSC-klagomål.docx", "A 50825-2025")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1261/klagomålsmail/A 50825-2025 FSC-klagomål mail.docx", "A 50825-2025")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1261/tillsyn/A 50825-2025 tillsynsbegäran.docx", "A 50825-2025")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1261/tillsynsmail/A 50825-2025 tillsynsbegäran mail.docx", "A 50825-2025")'

# Row 6
$ws.Range("A6").Value = "A 1468-2022"
$ws.Range("B6").Value = 44573
$ws.Range("C6").Value = 46064
$ws.Range("G6").Value = 1.8

# Row 7
$ws.Range("A7").Value = "A 65018-2023"
$ws.Range("B7").Value = 45287
$ws.Range("C7").Value = 46064
$ws.Range("G7").Value = 1.1

# Row 8
$ws.Range("C8").Value = 46064

# Row 9
$ws.Range("A9").Value = "A 53361-2024"
$ws.Range("B9").Value = 45614
$ws.Range("C9").Value = 46064
$ws.Range("G9").Value = 2.5

# Row 10
$ws.Range("A10").Value = "A 53361-2024"
$ws.Range("B10").Value = 45614
$ws.Range("C10").Value = 46064
$ws.Range("G10").Value = 0.4

# Row 11
$ws.Range("A11").Value = "A 18968-2025"
$ws.Range("B11").Value = 45764.53686342593
$ws.Range("C11").Value = 46064
$ws.Range("G11").Value = 8.199999999999999

# Row 12
$ws.Range("A12").Value = "A 53750-2025"
$ws.Range("B12").Value = 45960.65806712963
$ws.Range("C12").Value = 46064
$ws.Range("G12").Value = 0.9

# Row 13
$ws.Range("A13").Value = "A 19003-2025"
$ws.Range("B13").Value = 45764
$ws.Range("C13").Value = 46064
$ws.Range("G13").Value = 5.4
